$wb = $excel.ActiveWorkbook

# --- Sheet 1: "optimal models" ---
$ws1 = $wb.Worksheets.Item("optimal models")

$ws1.Range("A7").Value = "FRI"
$ws1.Range("B7").Value = 2
$ws1.Range("C7").Value = 0.75
$ws1.Range("D7").Value = 0.01
$ws1.Range("E7").Value = 2850
$ws1.Range("F7").Value = 0.9551
$ws1.Range("G7").Value = 0.90272
$ws1.Range("H7").Value = 0.9027
$ws1.Range("I7").Value = 0.174483

$ws1.Range("H8").Select()

# --- Sheet 2: "relative importance" ---
$ws2 = $wb.Worksheets.Item("relative importance")

$ws2.Range("A7").Value = "FRI"
$ws2.Range("B7").Value = 3.300034
$ws2.Range("C7").Value = 7.100325
$ws2.Range("D7").Value = 1.303535
$ws2.Range("E7").Value = 9.302117000000001
$ws2.Range("F7").Value = 10.819298
$ws2.Range("G7").Value = 9.662936999999999
$ws2.Range("H7").Value = 1.159564
$ws2.Range("I7").Value = 2.067737
$ws2.Range("J7").Value = 3.103566
$ws2.Range("K7").Value = 2.602923
$ws2.Range("L7").Value = 5.774734
$ws2.Range("M7").Value = 4.952779
$ws2.Range("N7").Value = 8.971748
$ws2.Range("O7").Value = 13.624163
$ws2.Range("P7").Value = 3.945235
$ws2.Range("Q7").Value = 8.502551
$ws2.Range("R7").Value = 3.806754

$ws2.Range("A8").Select()

# Restore the originally active sheet/tab (sheet1 "optimal models") and its
# selection so that tabSelected + selection stay where the diff expects.
$ws1.Activate()
$ws1.Range("H8").Select()
